$wb = $excel.ActiveWorkbook

# --- Sheet "Daily": row 2 updates ---
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 2939.88
$daily.Range("H2").Value = 6212.58
$daily.Range("I2").Value = 730.66
$daily.Range("J2").Value = 758.2
$daily.Range("L2").Value = 758.2

# --- Sheet "Hourly": rows 9-19 updates ---
$hourly = $wb.Worksheets.Item("Hourly")

$hourly.Range("I9").Value = 55.03
$hourly.Range("K9").Value = 1.81
$hourly.Range("M9").Value = 1.81

$hourly.Range("H10").Value = 108.15
$hourly.Range("I10").Value = 433.1
$hourly.Range("J10").Value = 48.99
$hourly.Range("K10").Value = 28.08
$hourly.Range("M10").Value = 28.08

$hourly.Range("H11").Value = 249.29
$hourly.Range("I11").Value = 635.22
$hourly.Range("J11").Value = 72.33
$hourly.Range("K11").Value = 66.31
$hourly.Range("M11").Value = 66.31

$hourly.Range("H12").Value = 369.73
$hourly.Range("I12").Value = 732.02
$hourly.Range("J12").Value = 85.84999999999999
$hourly.Range("K12").Value = 95.81
$hourly.Range("M12").Value = 95.81

$hourly.Range("H13").Value = 450.48
$hourly.Range("I13").Value = 780.1
$hourly.Range("J13").Value = 93.31
$hourly.Range("K13").Value = 115.85
$hourly.Range("M13").Value = 115.85

$hourly.Range("H14").Value = 481.8
$hourly.Range("I14").Value = 796.36
$hourly.Range("J14").Value = 95.95
$hourly.Range("K14").Value = 123.98
$hourly.Range("M14").Value = 123.98

$hourly.Range("H15").Value = 460.2
$hourly.Range("I15").Value = 785.3200000000001
$hourly.Range("J15").Value = 94.13
$hourly.Range("K15").Value = 118.87
$hourly.Range("M15").Value = 118.87

$hourly.Range("H16").Value = 388.06
$hourly.Range("I16").Value = 743.89
$hourly.Range("J16").Value = 87.62
$hourly.Range("K16").Value = 99.16
$hourly.Range("M16").Value = 99.16

$hourly.Range("H17").Value = 273.92
$hourly.Range("I17").Value = 658.4299999999999
$hourly.Range("J17").Value = 75.40000000000001
$hourly.Range("K17").Value = 68.48
$hourly.Range("M17").Value = 68.48

$hourly.Range("H18").Value = 134.46
$hourly.Range("I18").Value = 484.79
$hourly.Range("K18").Value = 35.03
$hourly.Range("M18").Value = 35.03

$hourly.Range("I19").Value = 108.32
$hourly.Range("K19").Value = 4.83
$hourly.Range("M19").Value = 4.83
